$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1728506666666667
$ws.Range("H2").Value = 0.518552
$ws.Range("I2").Value = 0.0840503369699626
$ws.Range("J2").Value = 0.0840503369699626
$ws.Range("M2").Value = 0.034325
$ws.Range("N2").Value = 0.102975
$ws.Range("O2").Value = 0.004508979075184418
$ws.Range("P2").Value = 0.004508979075184418
$ws.Range("Q2").Value = 0.005933099133333334
$ws.Range("R2").Value = 0.0533978922
$ws.Range("S2").Value = 0.0003789812106597607
$ws.Range("T2").Value = 0.0003789812106597607

$ws.Range("G3").Value = 0.1728506666666667
$ws.Range("H3").Value = 0.518552
$ws.Range("I3").Value = 0.0840503369699626
$ws.Range("J3").Value = 0.0840503369699626
$ws.Range("O3").Value = 0.9142039036746329
$ws.Range("P3").Value = 0.9142039036746329
$ws.Range("Q3").Value = 1.202946897321778
$ws.Range("R3").Value = 10.826522075896
$ws.Range("S3").Value = 0.07683914616310812
$ws.Range("T3").Value = 0.07683914616310812

$ws.Range("G4").Value = 0.1728506666666667
$ws.Range("H4").Value = 0.518552
$ws.Range("I4").Value = 0.0840503369699626
$ws.Range("J4").Value = 0.0840503369699626
$ws.Range("M4").Value = 0.5818573333333333
$ws.Range("N4").Value = 1.745572
$ws.Range("O4").Value = 0.0764335772976724
$ws.Range("P4").Value = 0.0764335772976724
$ws.Range("Q4").Value = 0.1005744279715556
$ws.Range("R4").Value = 0.9051698517440001
$ws.Range("S4").Value = 0.006424267927689049
$ws.Range("T4").Value = 0.006424267927689049

$ws.Range("G5").Value = 0.1728506666666667
$ws.Range("H5").Value = 0.518552
$ws.Range("I5").Value = 0.0840503369699626
$ws.Range("J5").Value = 0.0840503369699626
$ws.Range("M5").Value = 0.036948
$ws.Range("N5").Value = 0.110844
$ws.Range("O5").Value = 0.004853539952510238
$ws.Range("P5").Value = 0.004853539952510237
$ws.Range("Q5").Value = 0.006386486432000001
$ws.Range("R5").Value = 0.057478377888
$ws.Range("S5").Value = 0.0004079416685056618
$ws.Range("T5").Value = 0.0004079416685056617

$ws.Range("I6").Value = 0.6650661694281633
$ws.Range("J6").Value = 0.6650661694281633
$ws.Range("M6").Value = 0.034325
$ws.Range("N6").Value = 0.102975
$ws.Range("O6").Value = 0.004508979075184418
$ws.Range("P6").Value = 0.004508979075184418
$ws.Range("Q6").Value = 0.04694690890833333
$ws.Range("R6").Value = 0.422522180175
$ws.Range("S6").Value = 0.002998769441564643
$ws.Range("T6").Value = 0.002998769441564643

$ws.Range("I7").Value = 0.6650661694281633
$ws.Range("J7").Value = 0.6650661694281633
$ws.Range("O7").Value = 0.9142039036746329
$ws.Range("P7").Value = 0.9142039036746329
$ws.Range("S7").Value = 0.6080060882931617
$ws.Range("T7").Value = 0.6080060882931617

$ws.Range("I8").Value = 0.6650661694281633
$ws.Range("J8").Value = 0.6650661694281633
$ws.Range("M8").Value = 0.5818573333333333
$ws.Range("N8").Value = 1.745572
$ws.Range("O8").Value = 0.0764335772976724
$ws.Range("P8").Value = 0.0764335772976724
$ws.Range("Q8").Value = 0.7958165542795556
$ws.Range("R8").Value = 7.162348988516
$ws.Range("S8").Value = 0.05083338646905441
$ws.Range("T8").Value = 0.05083338646905441

$ws.Range("I9").Value = 0.6650661694281633
$ws.Range("J9").Value = 0.6650661694281633
$ws.Range("M9").Value = 0.036948
$ws.Range("N9").Value = 0.110844
$ws.Range("O9").Value = 0.004853539952510238
$ws.Range("P9").Value = 0.004853539952510237
$ws.Range("Q9").Value = 0.050534432348
$ws.Range("R9").Value = 0.454809891132
$ws.Range("S9").Value = 0.003227925224382534
$ws.Range("T9").Value = 0.003227925224382533

$ws.Range("G10").Value = 0.5159453333333334
$ws.Range("H10").Value = 1.547836
$ws.Range("I10").Value = 0.2508834936018741
$ws.Range("J10").Value = 0.2508834936018741
$ws.Range("M10").Value = 0.034325
$ws.Range("N10").Value = 0.102975
$ws.Range("O10").Value = 0.004508979075184418
$ws.Range("P10").Value = 0.004508979075184418
$ws.Range("Q10").Value = 0.01770982356666667
$ws.Range("R10").Value = 0.1593884121
$ws.Range("S10").Value = 0.001131228422960014
$ws.Range("T10").Value = 0.001131228422960014

$ws.Range("G11").Value = 0.5159453333333334
$ws.Range("H11").Value = 1.547836
$ws.Range("I11").Value = 0.2508834936018741
$ws.Range("J11").Value = 0.2508834936018741
$ws.Range("O11").Value = 0.9142039036746329
$ws.Range("P11").Value = 0.9142039036746329
$ws.Range("Q11").Value = 3.590699705647556
$ws.Range("R11").Value = 32.316297350828
$ws.Range("S11").Value = 0.2293586692183631
$ws.Range("T11").Value = 0.2293586692183631

$ws.Range("G12").Value = 0.5159453333333334
$ws.Range("H12").Value = 1.547836
$ws.Range("I12").Value = 0.2508834936018741
$ws.Range("J12").Value = 0.2508834936018741
$ws.Range("M12").Value = 0.5818573333333333
$ws.Range("N12").Value = 1.745572
$ws.Range("O12").Value = 0.0764335772976724
$ws.Range("P12").Value = 0.0764335772976724
$ws.Range("Q12").Value = 0.3002065757991111
$ws.Range("R12").Value = 2.701859182192
$ws.Range("S12").Value = 0.01917592290092894
$ws.Range("T12").Value = 0.01917592290092894

$ws.Range("G13").Value = 0.5159453333333334
$ws.Range("H13").Value = 1.547836
$ws.Range("I13").Value = 0.2508834936018741
$ws.Range("J13").Value = 0.2508834936018741
$ws.Range("M13").Value = 0.036948
$ws.Range("N13").Value = 0.110844
$ws.Range("O13").Value = 0.004853539952510238
$ws.Range("P13").Value = 0.004853539952510237
$ws.Range("Q13").Value = 0.019063148176
$ws.Range("R13").Value = 0.171568333584
$ws.Range("S13").Value = 0.001217673059622043
$ws.Range("T13").Value = 0.001217673059622042
